# 2017HandinStatistics.xlsx - "last 6 months work"
# Adds three new hand-in snapshots (columns J, K, L) mirroring the existing
# weekly columns B:I, extends the Mean/Sum summary formulas to cover them,
# and moves the sheet selection to the newly active area.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Clone formatting (date style etc.) from column I into J, K, L ------
# Row 3 is skipped on purpose: it only carries data in C:G (it is a merged
# banner/comment row) so J3/K3/L3 must stay empty, not inherit the row's
# default style.
$ws.Range("I2").Copy($ws.Range("J2"))
$ws.Range("I4:I16").Copy($ws.Range("J4:J16"))
$ws.Range("I17").Copy($ws.Range("J17"))
$ws.Range("I18").Copy($ws.Range("J18"))

$ws.Range("J2").Copy($ws.Range("K2"))
$ws.Range("J4:J16").Copy($ws.Range("K4:K16"))
$ws.Range("J17").Copy($ws.Range("K17"))
$ws.Range("J18").Copy($ws.Range("K18"))

$ws.Range("K2").Copy($ws.Range("L2"))
$ws.Range("K4:K16").Copy($ws.Range("L4:L16"))
$ws.Range("K17").Copy($ws.Range("L17"))
$ws.Range("K18").Copy($ws.Range("L18"))

# --- 2. New snapshot dates (row 2) -----------------------------------------
$ws.Cells.Item(2, 10).Value = 43122   # J2
$ws.Cells.Item(2, 11).Value = 43126   # K2
$ws.Cells.Item(2, 12).Value = 43136   # L2

# --- 3. New distribution counts (rows 4-16) ---------------------------------
$jVals = @(14, 7, 6, 8, 4, 4, 3, 6, 8, 11, 49, 94, 45)
$kVals = @(14, 7, 6, 8, 4, 4, 3, 6, 6, 4, 28, 112, 57)
$lVals = @(15, 7, 6, 6, 5, 5, 1, 4, 4, 2, 4, 130, 70)

for ($i = 0; $i -lt 13; $i++) {
    $row = 4 + $i
    $ws.Cells.Item($row, 10).Value = $jVals[$i]
    $ws.Cells.Item($row, 11).Value = $kVals[$i]
    $ws.Cells.Item($row, 12).Value = $lVals[$i]
}

# --- 4. Extend the "Mean" row (17) and "Sum" row (18) formulas -------------
$ws.Range("J17:K17").Formula = "=SUMPRODUCT(`$B`$4:`$B`$16,J4:J16)/SUM(J4:J16)"
$ws.Range("L17").Formula = "=SUMPRODUCT(`$B`$4:`$B`$16,L4:L16)/SUM(L4:L16)"

$ws.Range("J18:K18").Formula = "=SUM(J4:J16)"
$ws.Range("L18").Formula = "=SUM(L4:L16)"

# --- 5. Update the sheet view: scroll over one more column and move the ----
#        selection onto the newest snapshot column.
$ws.Range("L3").Select()

$wb.Save()
